$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$months = @{
    "Jan" = 1; "Feb" = 2; "Mar" = 3; "Apr" = 4; "May" = 5; "Jun" = 6;
    "Jul" = 7; "Aug" = 8; "Sep" = 9; "Oct" = 10; "Nov" = 11; "Dec" = 12
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $s = $cell.Value2

    if ($s -match '(\w+) (\d+) (\d+) (\d+):(\d+):(\d+) GMT-0500') {
        $mon = $months[$matches[1]]
        $day = [int]$matches[2]
        $year = [int]$matches[3]
        $hh = [int]$matches[4]
        $mm = [int]$matches[5]
        $ss = [int]$matches[6]

        $d = Get-Date -Year $year -Month $mon -Day $day -Hour $hh -Minute $mm -Second $ss
        $utc = $d.AddHours(5)
        $formatted = $utc.ToString("yyyy-MM-ddTHH:mm:ss") + ".000Z"

        $cell.Value2 = $formatted
    }
}
